$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark bugs #17, #18, #19 (rows 18-20) as Fixed
$ws.Range("F18").Value = "Fixed"
$ws.Range("F19").Value = "Fixed"
$ws.Range("F20").Value = "Fixed"

# Add new bug entry on row 22 (Id 21)
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Backend"
$ws.Range("C22").Value = "User"
$ws.Range("D22").Value = 'Add col "Loại" but can''t load the FKNavigation data'
$ws.Range("E22").Value = "Undone function"
$ws.Range("F22").Value = "Not fix"
$ws.Range("G22").Value = "Cá"

# Update the view: scroll so row 10 is the top-left visible row, and
# move the selection to G22
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("G22").Select()
